# GitHub Actions symbol-list refresh (cryptos.xlsx) on 2022-12-24.
# - Column D "Price" gets refreshed quotes for several coins.
# - Rows 9-17 ("One" .. "CoinExToken") shift down one rank: "One" moves
#   from rank 16 (row 17) up to rank 8 (row 9), pushing WazirX..CoinExToken
#   down by one row each; their rank prefix in column E is renumbered to
#   match the new row, and column D gets the refreshed price for the coin
#   now occupying that row.
# - Row 44 (LocalTraders) loses its "Bestin24h" suffix in column E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells (Coin name, Link, Volume label) - column B, C, E
$textCells = @{
  'B9' = 'One'
  'C9' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
  'E9' = '8OneONEBestin24h'
  'B10' = 'WazirX'
  'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
  'E10' = '9WazirXWRX'
  'B11' = 'MandalaExchangeToken'
  'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
  'E11' = '10MandalaExchangeTokenMDX'
  'B12' = 'LiechtensteinCryptoassetsExchange'
  'C12' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
  'E12' = '11LiechtensteinCryptoassetsExchangeLCX'
  'B13' = 'BitrueCoin'
  'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
  'E13' = '12BitrueCoinBTR'
  'B14' = 'BitMartToken'
  'C14' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
  'E14' = '13BitMartTokenBMX'
  'B15' = 'MCDex'
  'C15' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
  'E15' = '14MCDexMCB'
  'B16' = 'BitForexToken'
  'C16' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
  'E16' = '15BitForexTokenBF'
  'B17' = 'CoinExToken'
  'C17' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
  'E17' = '16CoinExTokenCET'
  'E44' = '43LocalTradersLCT'
}
foreach ($addr in $textCells.Keys) {
  $ws.Range($addr).Value = $textCells[$addr]
}

# Numeric-looking price cells stored as text (column D). The sheet keeps
# these as text (not numbers) so exact formatting - trailing/leading zeros,
# fixed decimal places, etc. - is preserved just like the source data.
# Setting NumberFormat to "@" (Text) before assigning the value stops Excel
# from auto-converting the numeric-looking string into a floating point
# number; resetting the style back to "Normal" afterwards drops the
# temporary text format so no stray formatting is left behind.
$priceCells = @{
  'D2' = '244.52'
  'D3' = '21.75'
  'D4' = '5.389'
  'D6' = '3.391'
  'D7' = '0.8154'
  'D8' = '0.9483'
  'D9' = '0.01120'
  'D10' = '0.1434'
  'D11' = '0.07433'
  'D12' = '0.03413'
  'D13' = '0.03053'
  'D14' = '0.09412'
  'D15' = '4.003'
  'D16' = '0.001592'
  'D17' = '0.04809'
  'D18' = '0.005516'
  'D19' = '0.004163'
  'D20' = '0.0009868'
  'D22' = '6.420'
  'D23' = '2.191'
  'D26' = '0.00007002'
  'D40' = '0.04017'
  'D41' = '0.006502'
  'D42' = '0.1074'
  'D44' = '0.006565'
  'D45' = '0.00005252'
  'D48' = '0.003108'
  'D49' = '0.00002100'
}
foreach ($addr in $priceCells.Keys) {
  $cell = $ws.Range($addr)
  $cell.NumberFormat = "@"
  $cell.Value = $priceCells[$addr]
  $cell.Style = "Normal"
}